# Algs.xlsx - add algorithm outline text and a small "Cell linked-list" diagram
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width (matches target stored width 40.5) ---
$ws.Columns.Item(2).ColumnWidth = 39.666666666666664

# --- Main outline text (column A/B), in order so shared strings line up ---
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "Цикл по всем ячейкам"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Для каждой поиск решений в 4-е стороны"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Для каждого решения поиск решений"

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Из всех цепочек решений поиск самых длинных"

# --- Small "Cell" chain diagram (E4:F8) ---
# Build the final cell style (thin border on all sides, centered) once on E5,
# then propagate it via copy/paste-special so no unused intermediate styles
# get left behind in styles.xml.
$e5 = $ws.Range("E5")
$e5.Borders.LineStyle = 1
$e5.HorizontalAlignment = -4108
$e5.VerticalAlignment = -4108

$e5.Copy()
$ws.Range("F5:F8").PasteSpecial(-4122)
$ws.Range("E6:E8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header cell above the chain: centered only, no border
$e4 = $ws.Range("E4")
$e4.Value = "Cell"
$e4.HorizontalAlignment = -4108
$e4.VerticalAlignment = -4108

# Chain values
$e5.Value = 1
$ws.Range("F5").Value = "Cell"

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = "null"

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = "null"

$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "Cell"

# --- Pseudo-code snippet (E12:E14) ---
$ws.Range("E14").Value = "ArrayList<Cell>"
$ws.Range("E13").Value = "Cell parentCell;"
$ws.Range("E12").Value = "Cell cell;"

# This string is first referenced last, so it gets the final shared-string index
$ws.Range("B7").Value = "Итог - arraylist из arraylist <Cells> и какой самый длинный тот и нужен"

# --- Selection matches the saved view ---
$null = $ws.Range("B8").Select()

Write-Output "done"
